# Generate Report for Handback
# Updates the handback-status report timestamps/status after regenerating
# the localization report (ht -> mt, and refreshed generate/handoff/handback
# datetimes for the zh-cn and de-de rows).

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest HO Xliff Generate Date" column (G) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-16 12:15:41"
$wsOverview.Range("G3").Value = "2016-08-16 12:15:41"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
# Priority column (E): ht -> mt
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"
# Correspond Handoff Datetime column (H)
$wsZhCn.Range("H2").Value = "2016-08-16 12:15:35"
$wsZhCn.Range("H3").Value = "2016-08-16 12:15:35"
# Correspond Handback DateTime column (K)
$wsZhCn.Range("K2").Value = "2016-08-16 12:16:08"
$wsZhCn.Range("K3").Value = "2016-08-16 12:16:08"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
# Priority column (E): ht -> mt
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"
# Correspond Handoff Datetime column (H)
$wsDeDe.Range("H2").Value = "2016-08-16 12:15:41"
$wsDeDe.Range("H3").Value = "2016-08-16 12:15:41"
# Correspond Handback DateTime column (K)
$wsDeDe.Range("K2").Value = "2016-08-16 12:16:15"
$wsDeDe.Range("K3").Value = "2016-08-16 12:16:15"
